# Auto-generated edit script applying the Sargatanas_Profits diff
# Updates leve-profit calculation cells (H,I,J,K,L,M,N) across all 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H87").Value = 42799.168
$ws.Range("J87").Value = 42799.168
$ws.Range("L87").Value = 42799.168
$ws.Range("N87").Value = -45295.168
$ws.Range("H90").Value = 42799.168
$ws.Range("J90").Value = 42799.168
$ws.Range("L90").Value = 128397.504
$ws.Range("N90").Value = -140877.504
$ws.Range("H100").Value = 2428.2942
$ws.Range("I100").Value = 1834.8
$ws.Range("J100").Value = 3276.1428
$ws.Range("K100").Value = 1834.8
$ws.Range("L100").Value = 3276.1428
$ws.Range("M100").Value = -1293.8
$ws.Range("N100").Value = -4358.1428
$ws.Range("H101").Value = 1046.3334
$ws.Range("J101").Value = 2342.5
$ws.Range("L101").Value = 7027.5
$ws.Range("N101").Value = -10271.5
$ws.Range("H103").Value = 857.375
$ws.Range("I103").Value = 446.25
$ws.Range("J103").Value = 994.4167
$ws.Range("K103").Value = 1338.75
$ws.Range("L103").Value = 2983.2501
$ws.Range("M103").Value = -752.75
$ws.Range("N103").Value = -4155.2501
$ws.Range("H125").Value = 166669200
$ws.Range("I125").Value = 166669200
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 1500022800
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -1500020340
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 4003.2727
$ws.Range("I132").Value = 3560.111
$ws.Range("K132").Value = 10680.333
$ws.Range("M132").Value = -8150.332999999999
$ws.Range("H137").Value = 2339.2727
$ws.Range("J137").Value = 1696.9
$ws.Range("L137").Value = 5090.700000000001
$ws.Range("N137").Value = -10190.7
$ws.Range("H141").Value = 2618.4
$ws.Range("I141").Value = 2273.25
$ws.Range("K141").Value = 6819.75
$ws.Range("M141").Value = -1639.75

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H4").Value = 125
$ws.Range("I4").Value = 125
$ws.Range("K4").Value = 125
$ws.Range("M4").Value = -9
$ws.Range("H45").Value = 4364.357
$ws.Range("I45").Value = 1266.2222
$ws.Range("K45").Value = 1266.2222
$ws.Range("M45").Value = -889.2221999999999
$ws.Range("H57").Value = 4552.4287
$ws.Range("I57").Value = 4552.4287
$ws.Range("K57").Value = 4552.4287
$ws.Range("M57").Value = -4068.4287
$ws.Range("H97").Value = 11928662
$ws.Range("I97").Value = 1821
$ws.Range("K97").Value = 1821
$ws.Range("M97").Value = -1325
$ws.Range("H110").Value = 12821428
$ws.Range("I110").Value = 779.1429000000001
$ws.Range("K110").Value = 779.1429000000001
$ws.Range("M110").Value = 1265.8571
$ws.Range("H126").Value = 5240.143
$ws.Range("I126").Value = 5240.143
$ws.Range("K126").Value = 15720.429
$ws.Range("M126").Value = -13250.429
$ws.Range("H132").Value = 5355.7734
$ws.Range("I132").Value = 4194.4614
$ws.Range("J132").Value = 8590.857
$ws.Range("K132").Value = 12583.3842
$ws.Range("L132").Value = 25772.571
$ws.Range("M132").Value = -10053.3842
$ws.Range("N132").Value = -30832.571

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H81").Value = 78569.71000000001
$ws.Range("J81").Value = 78569.71000000001
$ws.Range("L81").Value = 78569.71000000001
$ws.Range("N81").Value = -80691.71000000001
$ws.Range("H84").Value = 78569.71000000001
$ws.Range("J84").Value = 78569.71000000001
$ws.Range("L84").Value = 235709.13
$ws.Range("N84").Value = -246317.13
$ws.Range("H86").Value = 26598944
$ws.Range("I86").Value = 10872080
$ws.Range("K86").Value = 10872080
$ws.Range("M86").Value = -10870957
$ws.Range("H89").Value = 26598944
$ws.Range("I89").Value = 10872080
$ws.Range("K89").Value = 54360400
$ws.Range("M89").Value = -54354784
$ws.Range("H92").Value = 27000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 27000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 27000
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -31992
$ws.Range("H105").Value = 52135.87
$ws.Range("I105").Value = 75368.19
$ws.Range("J105").Value = 3348
$ws.Range("K105").Value = 75368.19
$ws.Range("L105").Value = 3348
$ws.Range("M105").Value = -73621.19
$ws.Range("N105").Value = -6842
$ws.Range("H113").Value = 5028.9
$ws.Range("I113").Value = 5028.9
$ws.Range("K113").Value = 5028.9
$ws.Range("M113").Value = -2858.9

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H22").Value = 283.33334
$ws.Range("I22").Value = 275
$ws.Range("K22").Value = 275
$ws.Range("M22").Value = 75
$ws.Range("H76").Value = 4956
$ws.Range("I76").Value = 4956
$ws.Range("K76").Value = 4956
$ws.Range("M76").Value = -4641
$ws.Range("H79").Value = 4956
$ws.Range("I79").Value = 4956
$ws.Range("K79").Value = 4956
$ws.Range("M79").Value = -3864
$ws.Range("H94").Value = 781.2593000000001
$ws.Range("I94").Value = 878.6667
$ws.Range("J94").Value = 703.3333
$ws.Range("K94").Value = 878.6667
$ws.Range("L94").Value = 703.3333
$ws.Range("M94").Value = -427.6667
$ws.Range("N94").Value = -1605.3333
$ws.Range("H132").Value = 6812.893
$ws.Range("I132").Value = 2253.2
$ws.Range("J132").Value = 9346.056
$ws.Range("K132").Value = 6759.599999999999
$ws.Range("L132").Value = 28038.168
$ws.Range("M132").Value = -4229.599999999999
$ws.Range("N132").Value = -33098.16800000001
$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H56").Value = 6880.4
$ws.Range("I56").Value = 6880.4
$ws.Range("K56").Value = 6880.4
$ws.Range("M56").Value = -6350.4
$ws.Range("H60").Value = 85.71429000000001
$ws.Range("J60").Value = 90
$ws.Range("L60").Value = 270
$ws.Range("N60").Value = -772

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H80").Value = 86823.586
$ws.Range("I80").Value = 2632.6667
$ws.Range("J80").Value = 171014.5
$ws.Range("K80").Value = 2632.6667
$ws.Range("L80").Value = 171014.5
$ws.Range("M80").Value = -1634.6667
$ws.Range("N80").Value = -173010.5
$ws.Range("H83").Value = 86823.586
$ws.Range("I83").Value = 2632.6667
$ws.Range("J83").Value = 171014.5
$ws.Range("K83").Value = 13163.3335
$ws.Range("L83").Value = 855072.5
$ws.Range("M83").Value = -8171.333500000001
$ws.Range("N83").Value = -865056.5
$ws.Range("H97").Value = 292.45456
$ws.Range("I97").Value = 292.45456
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 292.45456
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 203.54544
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 2665.5264
$ws.Range("I102").Value = 2605.7058
$ws.Range("K102").Value = 2605.7058
$ws.Range("M102").Value = -983.7058000000002
$ws.Range("H122").Value = 5954531.5
$ws.Range("I122").Value = 10205119
$ws.Range("J122").Value = 3709.2
$ws.Range("K122").Value = 30615357
$ws.Range("L122").Value = 11127.6
$ws.Range("M122").Value = -30612907
$ws.Range("N122").Value = -16027.6
$ws.Range("H126").Value = 2913.3333
$ws.Range("I126").Value = 2958.7856
$ws.Range("J126").Value = 2754.25
$ws.Range("K126").Value = 8876.356800000001
$ws.Range("L126").Value = 8262.75
$ws.Range("M126").Value = -6406.356800000001
$ws.Range("N126").Value = -13202.75
$ws.Range("H132").Value = 6885.4707
$ws.Range("I132").Value = 1785
$ws.Range("K132").Value = 5355
$ws.Range("M132").Value = -2825

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 7555.5
$ws.Range("I22").Value = 4000
$ws.Range("K22").Value = 4000
$ws.Range("M22").Value = -3705
$ws.Range("H27").Value = 7555.5
$ws.Range("I27").Value = 4000
$ws.Range("K27").Value = 4000
$ws.Range("M27").Value = -3893
$ws.Range("H40").Value = 5038.7085
$ws.Range("I40").Value = 3508.8572
$ws.Range("K40").Value = 3508.8572
$ws.Range("M40").Value = -3372.8572
$ws.Range("H93").Value = 7884.857
$ws.Range("I93").Value = 6566.1113
$ws.Range("K93").Value = 6566.1113
$ws.Range("M93").Value = -5318.1113

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H2").Value = 29999.5
$ws.Range("I2").Value = 29999.5
$ws.Range("K2").Value = 29999.5
$ws.Range("M2").Value = -29887.5
$ws.Range("H126").Value = 3587.9375
$ws.Range("I126").Value = 2474.5557
$ws.Range("J126").Value = 5019.4287
$ws.Range("K126").Value = 7423.6671
$ws.Range("L126").Value = 15058.2861
$ws.Range("M126").Value = -4953.6671
$ws.Range("N126").Value = -19998.2861
$ws.Range("H132").Value = 13900313
$ws.Range("I132").Value = 16671086
$ws.Range("J132").Value = 46450
$ws.Range("K132").Value = 50013258
$ws.Range("L132").Value = 139350
$ws.Range("M132").Value = -50010728
$ws.Range("N132").Value = -144410
